# Add a new text box ("TextBox 6") to slide 1 containing the project URL,
# auto-linked as a hyperlink - matching the authored edit:
#   <p:sp> id="7" name="TextBox 6" with text "https://github.com/platapp/plat-mvp"
#
# Note: the live editing session first created (and then undid) an earlier
# text box, which is why the surviving shape ends up with id=7 rather than
# id=6 - PowerPoint's per-slide shape-id counter keeps advancing across
# undo. We replicate that by adding+deleting a throwaway textbox first so
# the id counter lands on the same value.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Throwaway shape to reproduce the id counter advancing past the undone add.
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$placeholder.Delete()

# EMU -> points (1 pt = 12700 EMU) for the final shape's exact frame.
$left   = 1512147 / 12700
$top    = 5774803 / 12700
$width  = 6096000 / 12700
$height = 261610 / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "TextBox 6"
$shp.Fill.Visible = 0

$tf = $shp.TextFrame
$tf.WordWrap = -1

$tr = $tf.TextRange
$tr.Text = "https://github.com/platapp/plat-mvp"
$tr.Font.Size = 10.5

$tf.AutoSize = 1
$shp.Height = $height

$tr.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/platapp/plat-mvp"
